$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.576
$ws.Range("B4").Value = 6.093000000000001
$ws.Range("A6").Value = -21.325
$ws.Range("C6").Value = -13.063
$ws.Range("A7").Value = -20.296
$ws.Range("C7").Value = -11.879
$ws.Range("A8").Value = -20.996
$ws.Range("B8").Value = 6.411
$ws.Range("C8").Value = -13.371
$ws.Range("B9").Value = 5.808999999999999
$ws.Range("C10").Value = -13.026
$ws.Range("B12").Value = 5.52
$ws.Range("C13").Value = -12.918
$ws.Range("C14").Value = -12.987
$ws.Range("A16").Value = -21.2
$ws.Range("C16").Value = -12.485
$ws.Range("B17").Value = 5.619999999999999
$ws.Range("B18").Value = 5.605
$ws.Range("B19").Value = 6.919
$ws.Range("A20").Value = -21.854
$ws.Range("B20").Value = 5.773000000000001
$ws.Range("A21").Value = -20.186
$ws.Range("B26").Value = 6.132
$ws.Range("A28").Value = -21.555
$ws.Range("A29").Value = -21.5
$ws.Range("A30").Value = -22.344
$ws.Range("C30").Value = -13.511
$ws.Range("B31").Value = 5.966000000000001
$ws.Range("A32").Value = -21.583
$ws.Range("C37").Value = -13.258
$ws.Range("B39").Value = 6.718999999999999
$ws.Range("A40").Value = -21.432
$ws.Range("B40").Value = 6.167
$ws.Range("C40").Value = -11.547
$ws.Range("B41").Value = 6.794
$ws.Range("B42").Value = 6.615
$ws.Range("B43").Value = 5.903
$ws.Range("C44").Value = -12.906
$ws.Range("A46").Value = -21.63
$ws.Range("B47").Value = 5.833
$ws.Range("B48").Value = 5.443
$ws.Range("A51").Value = -21.554
$ws.Range("A52").Value = -21.651
$ws.Range("B54").Value = 5.436000000000001
$ws.Range("A57").Value = -21.889
$ws.Range("A59").Value = -22.154
$ws.Range("A62").Value = -21.792
$ws.Range("B62").Value = 5.983000000000001
$ws.Range("B63").Value = 5.252
$ws.Range("B64").Value = 5.527
$ws.Range("A66").Value = -21.563
$ws.Range("C70").Value = -11.308
$ws.Range("A73").Value = -21.454
$ws.Range("A74").Value = -20.688
$ws.Range("B76").Value = 7.201000000000001
$ws.Range("A77").Value = -21.267
$ws.Range("B81").Value = 5.911
$ws.Range("B84").Value = 6.161
$ws.Range("B89").Value = 5.261
$ws.Range("C89").Value = -13.71
$ws.Range("C91").Value = -12.846
$ws.Range("A92").Value = -21.655
$ws.Range("C93").Value = -10.411
$ws.Range("B94").Value = 5.888000000000001
$ws.Range("C98").Value = -13.295
$ws.Range("A100").Value = -21.758
